$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_data3")

# Move the table header from C14 to B14 (same text)
$ws.Range("C14").Value = $null
$ws.Range("B14").Value = "Sweden_Pay_Now_Direct_debit"

# Replace the old column_name#/value# rows (C15:D17) with key4/value4 pairs in B15:C17
$ws.Range("C15").Value = $null
$ws.Range("D15").Value = $null
$ws.Range("B15").Value = "key4"
$ws.Range("C15").Value = "value4"

$ws.Range("C16").Value = $null
$ws.Range("D16").Value = $null
$ws.Range("B16").Value = "key4"
$ws.Range("C16").Value = "value4"

$ws.Range("C17").Value = $null
$ws.Range("D17").Value = $null
$ws.Range("B17").Value = "key4"
$ws.Range("C17").Value = "value4"

# Row 18 becomes the header of a second exported table
$ws.Range("C18").Value = $null
$ws.Range("D18").Value = $null
$ws.Range("B18").Value = "Sweden_Pay_Now_Card"
$ws.Range("C18").Value = "value4"

# New rows 19-22 with key4/value4 pairs for the second table
$ws.Range("B19").Value = "key4"
$ws.Range("C19").Value = "value4"

$ws.Range("B20").Value = "key4"
$ws.Range("C20").Value = "value4"

$ws.Range("B21").Value = "key4"
$ws.Range("C21").Value = "value4"

$ws.Range("B22").Value = "key4"
$ws.Range("C22").Value = "value4"
